$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.65
$ws.Range("I2").Value = 5.8
$ws.Range("K2").Value = 4.7
$ws.Range("N2").Value = 5
$ws.Range("O2").Value = 1.22
$ws.Range("Q2").Value = 1.67
$ws.Range("R2").Value = 1.56
$ws.Range("U2").Value = 2.24
$ws.Range("V2").Value = 1.2
$ws.Range("W2").Value = 2.52
$ws.Range("Y2").Value = 25
$ws.Range("AB2").Value = 11
$ws.Range("AL2").Value = 29
$ws.Range("AN2").Value = 7.2
$ws.Range("G3").Value = 2.5
$ws.Range("Q3").Value = 2.68
$ws.Range("V3").Value = 1.26
$ws.Range("F4").Value = 2.34
$ws.Range("G4").Value = 2.48
$ws.Range("H4").Value = 2.84
$ws.Range("I4").Value = 3.05
$ws.Range("J4").Value = 3.85
$ws.Range("K4").Value = 4.2
$ws.Range("Q4").Value = 1.53
$ws.Range("R4").Value = 1.62
$ws.Range("V4").Value = 1.48
$ws.Range("W4").Value = 1.67
$ws.Range("AB4").Value = 17
$ws.Range("AC4").Value = 10.5
$ws.Range("AD4").Value = 14
$ws.Range("AF4").Value = 21
$ws.Range("AH4").Value = 18
$ws.Range("AN4").Value = 13
$ws.Range("H5").Value = 3.1
$ws.Range("L5").Value = 1.55
$ws.Range("M5").Value = 1.15
$ws.Range("N5").Value = 2.2
$ws.Range("O5").Value = 1.64
$ws.Range("Q5").Value = 2.66
$ws.Range("S5").Value = 1.05
$ws.Range("T5").Value = 2.24
$ws.Range("W5").Value = 1.51
$ws.Range("X5").Value = 7.2
$ws.Range("Z5").Value = 26
$ws.Range("AA5").Value = 95
$ws.Range("AB5").Value = 8.800000000000001
$ws.Range("AE5").Value = 75
$ws.Range("AH5").Value = 34
$ws.Range("AJ5").Value = 60
$ws.Range("AK5").Value = 60
$ws.Range("F6").Value = 1.52
$ws.Range("G6").Value = 1.64
$ws.Range("H6").Value = 6.2
$ws.Range("I6").Value = 8.800000000000001
$ws.Range("J6").Value = 4
$ws.Range("K6").Value = 5.1
$ws.Range("M6").Value = 1.05
$ws.Range("N6").Value = 4.2
$ws.Range("P6").Value = 2.14
$ws.Range("Q6").Value = 1.7
$ws.Range("R6").Value = 1.44
$ws.Range("S6").Value = 2.78
$ws.Range("T6").Value = 1.82
$ws.Range("U6").Value = 1.98
$ws.Range("V6").Value = 1.14
$ws.Range("W6").Value = 2.56
$ws.Range("X6").Value = 24
$ws.Range("Y6").Value = 26
$ws.Range("AA6").Value = 230
$ws.Range("AB6").Value = 11.5
$ws.Range("AC6").Value = 12.5
$ws.Range("AD6").Value = 34
$ws.Range("AE6").Value = 120
$ws.Range("AF6").Value = 12
$ws.Range("AG6").Value = 12.5
$ws.Range("AH6").Value = 27
$ws.Range("AI6").Value = 110
$ws.Range("AJ6").Value = 18
$ws.Range("AK6").Value = 20
$ws.Range("AL6").Value = 40
$ws.Range("AM6").Value = 140
$ws.Range("AN6").Value = 9.199999999999999
$ws.Range("F7").Value = 2.48
$ws.Range("R7").Value = 1.1
$ws.Range("S7").Value = 2.12
$ws.Range("X7").Value = 17
$ws.Range("Y7").Value = 16
$ws.Range("G8").Value = 1.77
$ws.Range("N8").Value = 2
$ws.Range("Q8").Value = 2.42
$ws.Range("W8").Value = 2.3
$ws.Range("F9").Value = 1.67
$ws.Range("G9").Value = 1.69
$ws.Range("H10").Value = 27
$ws.Range("J10").Value = 11
$ws.Range("P10").Value = 2.9
$ws.Range("R10").Value = 1.76
$ws.Range("S10").Value = 2.08
$ws.Range("T10").Value = 2.68
$ws.Range("U10").Value = 1.47
$ws.Range("X10").Value = 990
$ws.Range("AF10").Value = 7.4
